# Commit: "fix demo data sheet name"
#
# Rename the three "Step" sheets to the "処理" (processing) naming so they
# read consistently with the rest of the workbook tabs.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("整形Step1").Name = "整形処理1"
$wb.Worksheets.Item("整形Step2").Name = "整形処理2"
$wb.Worksheets.Item("整形Step3").Name = "整形処理3"

# The author had the "初期データ" sheet active/selected when the file was
# last saved (rather than "元データ") -- bring that tab to the front to
# match.
$wb.Worksheets.Item("初期データ").Activate()
